$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21, shifting existing rows 21-34 down to 22-35
$ws.Rows("21:21").Insert()

# Populate the newly inserted row 21 with the new data record
$ws.Cells.Item(21, 1).Value = 8
$ws.Cells.Item(21, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(21, 3).Value = "Coquimbo"
$ws.Cells.Item(21, 4).Value = 45134
$ws.Cells.Item(21, 5).Value = 4
$ws.Cells.Item(21, 6).Value = 100112013
$ws.Cells.Item(21, 7).Value = "Alcachofa"
$ws.Cells.Item(21, 8).Value = "Española"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 500
$ws.Cells.Item(21, 11).Value = 12500
$ws.Cells.Item(21, 12).Value = 13000
$ws.Cells.Item(21, 13).Value = 12750
$ws.Cells.Item(21, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(21, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(21, 16).Value = 425
$ws.Cells.Item(21, 17).Value = 30
$ws.Cells.Item(21, 18).Value = "Hortaliza"
